# Apply updated cryptos list values (scraped prices / 1h volume changes)
# and the row-41/42 coin swap (Aave <-> FraxShare), per commit:
# "Updated cryptos list on Sun Oct  1 14:41:33 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.187.83'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '1.687.38'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '215.56'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '0.520'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '23.17'
$ws.Range('E8').Value = '  +10.23%  '
$ws.Range('D9').Value = '0.261'
$ws.Range('E9').Value = '  +4.11%  '
$ws.Range('D10').Value = '0.0627'
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('D11').Value = '0.0890'
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('D12').Value = '1.929.37'
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('D13').Value = '1.690.32'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').Value = '4.20'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').Value = '0.555'
$ws.Range('E15').Value = '  +4.49%  '
$ws.Range('D16').Value = '67.16'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').Value = '27.226.31'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').Value = '236.62'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '8.07'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('D20').Value = '0.0₃0744'
$ws.Range('E20').Value = '  +1.24%  '
$ws.Range('D22').Value = '4.55'
$ws.Range('E22').Value = '  +2.72%  '
$ws.Range('D23').Value = '9.61'
$ws.Range('E23').Value = '  +4.43%  '
$ws.Range('D24').Value = '2.07'
$ws.Range('E24').Value = '  -2.51%  '
$ws.Range('D25').Value = '147.38'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('D26').Value = '7.33'
$ws.Range('E26').Value = '  +1.64%  '
$ws.Range('D27').Value = '16.47'
$ws.Range('E27').Value = '  +2.59%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '0.0505'
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').Value = '3.39'
$ws.Range('E32').Value = '  +1.93%  '
$ws.Range('D33').Value = '1.543.79'
$ws.Range('E33').Value = '  +3.14%  '
$ws.Range('E34').Value = '  +2.18%  '
$ws.Range('D35').Value = '1.66'
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('D36').Value = '0.606'
$ws.Range('E36').Value = '  +3.47%  '
$ws.Range('D37').Value = '0.946'
$ws.Range('E37').Value = '  +3.49%  '
$ws.Range('D38').Value = '2.39'
$ws.Range('E38').Value = '  -0.59%  '
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '5.76'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '69.31'
$ws.Range('E42').Value = '  +2.64%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = '2.26'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').Value = '1.836.57'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').Value = '0.794'
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('D47').Value = '90.58'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('E48').Value = '  +5.43%  '
$ws.Range('D49').Value = '1.62'
$ws.Range('E49').Value = '  +6.00%  '
$ws.Range('D50').Value = '8.37'
$ws.Range('E50').Value = '  +7.32%  '
$ws.Range('D51').Value = '0.104'
$ws.Range('E51').Value = '  +0.73%  '
